# NotesData.xlsx — "added mari notes speed, remove old lyrics files"
#
# Fills in the "speed" column (D) on the three "mari" sheets (previously all
# placeholder 1's, two rows per sheet still blank) with the tuned per-note
# speed values, and restores the UI selection/active-sheet state left behind
# after editing.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# mari1 — speed column (D2:D19)
# ---------------------------------------------------------------------
$mari1 = $wb.Worksheets.Item("mari1")
$mari1Speeds = @{
    2  = 3
    3  = 3
    4  = 1.5
    5  = 0.5
    6  = 1.2
    7  = 1.2
    8  = 0.8
    9  = 1.3
    10 = 0.5
    11 = 1.6
    12 = 1.1
    13 = 1.4
    14 = 0.5
    15 = 1.5
    16 = 0.8
    17 = 1.7
    18 = 0.3
    19 = 1.5
}
foreach ($row in $mari1Speeds.Keys) {
    $mari1.Range("D$row").Value = $mari1Speeds[$row]
}

# ---------------------------------------------------------------------
# mari2 — speed column (D2:D19); D14, D16, D17 stay at their old value
# ---------------------------------------------------------------------
$mari2 = $wb.Worksheets.Item("mari2")
$mari2Speeds = @{
    2  = 3
    3  = 3
    4  = 1.5
    5  = 1.6
    6  = 0.5
    7  = 1.3
    8  = 0.9
    9  = 1.6
    10 = 0.3
    11 = 1.3
    12 = 1.1
    13 = 1.1
    15 = 0.9
    18 = 0.4
    19 = 1.4
}
foreach ($row in $mari2Speeds.Keys) {
    $mari2.Range("D$row").Value = $mari2Speeds[$row]
}

# ---------------------------------------------------------------------
# Selections left on each sheet
# ---------------------------------------------------------------------
$mari1.Range("D20").Select()

$mari3 = $wb.Worksheets.Item("mari3")

# mari2 ends up the active tab, with its selection parked at D20
$mari2.Activate()
$mari2.Range("D20").Select()

# mari3 keeps its own last selection (C16) but loses the "active" flag,
# which Activate() above already took care of.
$mari3.Range("C16").Select()

# Re-affirm mari2 as the active sheet/tab (matches activeTab 8 -> 7).
$mari2.Activate()
